$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 8000.5
$ws.Range("I43").Value = 8000.5
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 8000.5
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -7931.5
$ws.Range("N43").ClearContents()
$ws.Range("H87").Value = 79998.5
$ws.Range("J87").Value = 79998.5
$ws.Range("L87").Value = 79998.5
$ws.Range("N87").Value = -82494.5
$ws.Range("H88").Value = 6244.5
$ws.Range("I88").Value = 4990
$ws.Range("K88").Value = 4990
$ws.Range("M88").Value = -4584
$ws.Range("H90").Value = 79998.5
$ws.Range("J90").Value = 79998.5
$ws.Range("L90").Value = 239995.5
$ws.Range("N90").Value = -252475.5
$ws.Range("H91").Value = 6244.5
$ws.Range("I91").Value = 4990
$ws.Range("K91").Value = 4990
$ws.Range("M91").Value = -3586
$ws.Range("H103").Value = 567.1667
$ws.Range("I103").Value = 560.75
$ws.Range("J103").Value = 580
$ws.Range("K103").Value = 1682.25
$ws.Range("L103").Value = 1740
$ws.Range("M103").Value = -1096.25
$ws.Range("N103").Value = -2912
$ws.Range("H107").Value = 1206.9166
$ws.Range("I107").Value = 976.44446
$ws.Range("J107").Value = 1898.3334
$ws.Range("K107").Value = 976.44446
$ws.Range("L107").Value = 1898.3334
$ws.Range("M107").Value = 943.55554
$ws.Range("N107").Value = -5738.3334
$ws.Range("H132").Value = 3991
$ws.Range("I132").Value = 3991
$ws.Range("K132").Value = 11973
$ws.Range("M132").Value = -9443
$ws.Range("H135").Value = 744.9286
$ws.Range("I135").Value = 744.9286
$ws.Range("K135").Value = 6704.3574
$ws.Range("M135").Value = -4169.3574

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4428.9644
$ws.Range("I32").Value = 3427.3845
$ws.Range("K32").Value = 3427.3845
$ws.Range("M32").Value = -3140.3845
$ws.Range("H46").Value = 12405.286
$ws.Range("I46").Value = 15379
$ws.Range("J46").Value = 10175
$ws.Range("K46").Value = 15379
$ws.Range("L46").Value = 10175
$ws.Range("M46").Value = -15060
$ws.Range("N46").Value = -10813
$ws.Range("H74").Value = 1607.3529
$ws.Range("I74").Value = 1551.8572
$ws.Range("J74").Value = 1866.3334
$ws.Range("K74").Value = 1551.8572
$ws.Range("L74").Value = 1866.3334
$ws.Range("M74").Value = -677.8571999999999
$ws.Range("N74").Value = -3614.3334
$ws.Range("H77").Value = 1607.3529
$ws.Range("I77").Value = 1551.8572
$ws.Range("J77").Value = 1866.3334
$ws.Range("K77").Value = 7759.286
$ws.Range("L77").Value = 9331.666999999999
$ws.Range("M77").Value = -3391.286
$ws.Range("N77").Value = -18067.667
$ws.Range("H122").Value = 1885.2
$ws.Range("I122").Value = 1885.2
$ws.Range("K122").Value = 5655.6
$ws.Range("M122").Value = -3205.6

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3736.5
$ws.Range("I20").Value = 2918.75
$ws.Range("K20").Value = 2918.75
$ws.Range("M20").Value = -2671.75
$ws.Range("H107").Value = 1064.75
$ws.Range("I107").Value = 1064.75
$ws.Range("K107").Value = 1064.75
$ws.Range("M107").Value = 855.25
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 606
$ws.Range("I2").Value = 606
$ws.Range("K2").Value = 606
$ws.Range("M2").Value = -493
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H31").Value = 2139.5
$ws.Range("I31").Value = 2139.5
$ws.Range("K31").Value = 2139.5
$ws.Range("M31").Value = -1844.5
$ws.Range("H34").Value = 2139.5
$ws.Range("I34").Value = 2139.5
$ws.Range("K34").Value = 2139.5
$ws.Range("M34").Value = -1937.5
$ws.Range("H58").Value = 6179.8
$ws.Range("I58").Value = 5224.75
$ws.Range("K58").Value = 5224.75
$ws.Range("M58").Value = -5021.75
$ws.Range("H107").Value = 1138.8889
$ws.Range("I107").Value = 613.7143
$ws.Range("K107").Value = 613.7143
$ws.Range("M107").Value = 1306.2857
$ws.Range("H134").Value = 3378.8235
$ws.Range("I134").Value = 3502.5
$ws.Range("J134").Value = 1400
$ws.Range("K134").Value = 10507.5
$ws.Range("L134").Value = 4200
$ws.Range("M134").Value = -7972.5
$ws.Range("N134").Value = -9270
$ws.Range("H136").Value = 6179.8
$ws.Range("I136").Value = 5224.75
$ws.Range("K136").Value = 15674.25
$ws.Range("M136").Value = -13124.25
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1516099.5
$ws.Range("I4").Value = 21213.572
$ws.Range("K4").Value = 63640.716
$ws.Range("M4").Value = -63528.716
$ws.Range("H132").Value = 1242.5834
$ws.Range("I132").Value = 1166.1
$ws.Range("K132").Value = 10494.9
$ws.Range("M132").Value = -7964.9

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 7849.2856
$ws.Range("I57").Value = 2989.6
$ws.Range("J57").Value = 19998.5
$ws.Range("K57").Value = 2989.6
$ws.Range("L57").Value = 19998.5
$ws.Range("M57").Value = -2169.6
$ws.Range("N57").Value = -21638.5
$ws.Range("H70").Value = 18604.846
$ws.Range("I70").Value = 24607.334
$ws.Range("K70").Value = 24607.334
$ws.Range("M70").Value = -24337.334
$ws.Range("H73").Value = 18604.846
$ws.Range("I73").Value = 24607.334
$ws.Range("K73").Value = 24607.334
$ws.Range("M73").Value = -23671.334
$ws.Range("H97").Value = 638
$ws.Range("I97").Value = 614.8570999999999
$ws.Range("K97").Value = 614.8570999999999
$ws.Range("M97").Value = -118.8570999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8162.091
$ws.Range("I7").Value = 8088.5
$ws.Range("J7").Value = 8898
$ws.Range("K7").Value = 8088.5
$ws.Range("L7").Value = 8898
$ws.Range("M7").Value = -7976.5
$ws.Range("N7").Value = -9122
$ws.Range("H40").Value = 3316.6667
$ws.Range("I40").Value = 3316.6667
$ws.Range("K40").Value = 3316.6667
$ws.Range("M40").Value = -3180.6667
$ws.Range("H100").Value = 925.7143
$ws.Range("I100").Value = 830
$ws.Range("K100").Value = 830
$ws.Range("M100").Value = -289
$ws.Range("H126").Value = 8162.091
$ws.Range("I126").Value = 8088.5
$ws.Range("J126").Value = 8898
$ws.Range("K126").Value = 24265.5
$ws.Range("L126").Value = 26694
$ws.Range("M126").Value = -21795.5
$ws.Range("N126").Value = -31634

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H81").Value = 4308.6
$ws.Range("I81").Value = 2246.5
$ws.Range("K81").Value = 4493
$ws.Range("M81").Value = -3432
$ws.Range("H84").Value = 4308.6
$ws.Range("I84").Value = 2246.5
$ws.Range("K84").Value = 22465
$ws.Range("M84").Value = -17161
$ws.Range("H136").Value = 6781.1113
$ws.Range("I136").Value = 2819
$ws.Range("K136").Value = 8457
$ws.Range("M136").Value = -5907

Write-Output "Applied 190 cell changes across 8 sheets"